$d = $word.ActiveDocument

# Word "green" accent color used throughout this doc for highlighted bullets (RGB 00B050).
# OLE/Word color longs are 0x00BBGGRR, so RGB 00,B0,50 -> 0x0050B000.
$green = 5287936

# ---------------------------------------------------------------------------
# Locate the relevant paragraphs by their text rather than hard-coded indices
# so the script is resilient to minor shifts elsewhere in the document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count

$pStore = $null
$pBackup = $null
$pMysql = $null
$pOneTime = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($pStore -eq $null -and $t.StartsWith("Store images as blobs")) {
        $pStore = $p
    }
    elseif ($pStore -ne $null -and $pBackup -eq $null -and $t.StartsWith("Backup and restore")) {
        $pBackup = $p
    }
    elseif ($pBackup -ne $null -and $pMysql -eq $null -and $t.StartsWith("Migration to")) {
        $pMysql = $p
    }
    elseif ($pMysql -ne $null -and $pOneTime -eq $null -and $t.StartsWith("One time code for alarm manager")) {
        $pOneTime = $p
    }
}

# ---------------------------------------------------------------------------
# 1. Drop the old _GoBack bookmark that currently sits between the "d" and
#    "b" runs of "db" in the "Store images as blobs in db" bullet.
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2. "Store images as blobs in db (instead of files)..." bullet: paint it
#    green (paragraph mark + every run) and merge the now-adjacent "d"/"b"
#    runs into a single "db" run.
# ---------------------------------------------------------------------------
$storeText = $pStore.Range.Text
$dbOffset = $storeText.IndexOf("db")
$pStore.Range.Font.Color = $green

$dbRange = $d.Range($pStore.Range.Start + $dbOffset, $pStore.Range.Start + $dbOffset + 2)
$dbRange.Find.Execute("db", $true, $false, $false, $false, $false, $true, 1, $false, "db", 2)

# ---------------------------------------------------------------------------
# 3. "Backup and restore" bullet -> green.
# ---------------------------------------------------------------------------
$pBackup.Range.Font.Color = $green

# ---------------------------------------------------------------------------
# 4. "Migration to mysql" bullet -> green.
# ---------------------------------------------------------------------------
$pMysql.Range.Font.Color = $green

# ---------------------------------------------------------------------------
# 5. Re-insert the _GoBack bookmark, now collapsed at the very start of the
#    "One time code for alarm manager" bullet.
# ---------------------------------------------------------------------------
$goBackRange = $pOneTime.Range.Duplicate
$goBackRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
